# Auto-generated script to apply F/G column numeric updates per the diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 225
$ws.Range("F4").Value = 407
$ws.Range("F6").Value = 10389
$ws.Range("F8").Value = 3607
$ws.Range("F9").Value = 214
$ws.Range("F10").Value = 2464
$ws.Range("F11").Value = 46
$ws.Range("F12").Value = 2880
$ws.Range("F15").Value = 2217
$ws.Range("F17").Value = 102
$ws.Range("F19").Value = 405
$ws.Range("F21").Value = 165
$ws.Range("F22").Value = 324
$ws.Range("F23").Value = 280
$ws.Range("F24").Value = 256
$ws.Range("F26").Value = 1346
$ws.Range("F27").Value = 23
$ws.Range("F28").Value = 1269
$ws.Range("F29").Value = 110
$ws.Range("F30").Value = 137
$ws.Range("F32").Value = 3946
$ws.Range("F33").Value = 3402
$ws.Range("F34").Value = 45
$ws.Range("F36").Value = 1063
$ws.Range("F37").Value = 415
$ws.Range("F39").Value = 1297
$ws.Range("F42").Value = 78
$ws.Range("F44").Value = 46
$ws.Range("F45").Value = 22

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G4").Value = 280
$ws.Range("F15").Value = 39

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 760
$ws.Range("F5").Value = 2120

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 225
$ws.Range("F3").Value = 760
$ws.Range("F7").Value = 407
$ws.Range("F11").Value = 10389
$ws.Range("G12").Value = 280
$ws.Range("F14").Value = 3608
$ws.Range("F15").Value = 2464
$ws.Range("F16").Value = 46
$ws.Range("F17").Value = 2880
$ws.Range("F19").Value = 2217
$ws.Range("F21").Value = 102
$ws.Range("F23").Value = 405
$ws.Range("F25").Value = 324
$ws.Range("F26").Value = 280
$ws.Range("F27").Value = 256
$ws.Range("F29").Value = 1346
$ws.Range("F30").Value = 23
$ws.Range("F31").Value = 1269
$ws.Range("F32").Value = 110
$ws.Range("F33").Value = 137
$ws.Range("F36").Value = 3946
$ws.Range("F37").Value = 3402
$ws.Range("F38").Value = 45
$ws.Range("F39").Value = 1063
$ws.Range("F41").Value = 415
$ws.Range("F43").Value = 39
$ws.Range("F44").Value = 1297
$ws.Range("F46").Value = 78
$ws.Range("F47").Value = 46
$ws.Range("F48").Value = 22

